$wb = $excel.ActiveWorkbook

# Sheet 1
$ws = $wb.Worksheets.Item(1)
$ws.Range("F8").Value = 13
$ws.Range("F9").Value = 340
$ws.Range("F11").Value = 1267
$ws.Range("F12").Value = 29060
$ws.Range("F13").Value = 4125
$ws.Range("F14").Value = 42
$ws.Range("F15").Value = 262
$ws.Range("F16").Value = 484
$ws.Range("F17").Value = 38
$ws.Range("F19").Value = 12
$ws.Range("C20").Value = "广州·pokemon only PMO 2024-得闲饮茶"
$ws.Range("F21").Value = 340
$ws.Range("F22").Value = 629
$ws.Range("F24").Value = 281
$ws.Range("F27").Value = 69
$ws.Range("F29").Value = 665
$ws.Range("F30").Value = 214
$ws.Range("F35").Value = 638
$ws.Range("F38").Value = 8

# Sheet 2
$ws = $wb.Worksheets.Item(2)
$ws.Range("G4").Value = 102
$ws.Range("F6").Value = 385
$ws.Range("F7").Value = 885
$ws.Range("F10").Value = 90
$ws.Range("F23").Value = 4247

# Sheet 3
$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 262
$ws.Range("F4").Value = 1199

# Sheet 4
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 262
$ws.Range("F4").Value = 1199
$ws.Range("G6").Value = 102
$ws.Range("F7").Value = 385
$ws.Range("F9").Value = 885
$ws.Range("F14").Value = 13
$ws.Range("F15").Value = 340
$ws.Range("F18").Value = 1267
$ws.Range("F19").Value = 90
$ws.Range("F20").Value = 90
$ws.Range("F28").Value = 484
$ws.Range("F29").Value = 38
$ws.Range("F30").Value = 12
$ws.Range("C32").Value = "广州·pokemon only PMO 2024-得闲饮茶"
$ws.Range("F34").Value = 340
$ws.Range("F35").Value = 629
$ws.Range("F38").Value = 69
$ws.Range("F40").Value = 665
$ws.Range("F42").Value = 214
$ws.Range("F48").Value = 638
